# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets
# to reflect newly scraped data (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 2941
    $ws.Range("F7").Value = 1591
    $ws.Range("F11").Value = 1312
    $ws.Range("F13").Value = 431
    $ws.Range("F16").Value = 56
    $ws.Range("F20").Value = 2999
    $ws.Range("F21").Value = 365
    $ws.Range("F22").Value = 49
}
